$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.870.54"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "3.545.98"
$ws.Range("E3").Value = "  +4.37%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'599.01"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").Value = "'136.28"
$ws.Range("E6").Value = "  +3.60%  "
$ws.Range("D7").Value = "3.544.38"
$ws.Range("E7").Value = "  +4.34%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +3.66%  "
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("D11").Value = "'6.93"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'0.387"
$ws.Range("E12").Value = "  +4.54%  "
$ws.Range("D13").Value = "4.147.29"
$ws.Range("E13").Value = "  +4.32%  "
$ws.Range("E14").Value = "  +3.97%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'27.10"
$ws.Range("E15").Value = "  +5.09%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.551.08"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "64.797.89"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'10.04"
$ws.Range("E19").Value = "  +6.30%  "
$ws.Range("D20").Value = "'14.43"
$ws.Range("E20").Value = "  +7.50%  "
$ws.Range("D21").Value = "'5.84"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("D22").Value = "'389.07"
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("E23").Value = "  +7.03%  "
$ws.Range("D24").Value = "3.691.77"
$ws.Range("E24").Value = "  +4.48%  "
$ws.Range("D25").Value = "'74.15"
$ws.Range("E25").Value = "  +3.86%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  +13.88%  "
$ws.Range("D28").Value = "'7.66"
$ws.Range("E28").Value = "  +8.63%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +5.57%  "
$ws.Range("D31").Value = "'8.35"
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.46"
$ws.Range("E32").Value = "  +25.47%  "
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D33").Value = "3.555.73"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "'23.97"
$ws.Range("E34").Value = "  +5.26%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("D37").Value = "'170.18"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +5.73%  "
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("D40").Value = "'5.00"
$ws.Range("E40").Value = "  +10.31%  "
$ws.Range("D41").Value = "'0.0807"
$ws.Range("E41").Value = "  +7.45%  "
$ws.Range("E42").Value = "  +4.51%  "
$ws.Range("D43").Value = "'26.96"
$ws.Range("E43").Value = "  +22.97%  "
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'4.46"
$ws.Range("E46").Value = "  +5.48%  "
$ws.Range("D47").Value = "'1.21"
$ws.Range("E47").Value = "  +10.56%  "
$ws.Range("D48").Value = "'1.65"
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("D49").Value = "'6.91"
$ws.Range("E49").Value = "  +7.25%  "
$ws.Range("D50").Value = "2.443.28"
$ws.Range("E50").Value = "  +12.19%  "
$ws.Range("D51").Value = "'2.36"
$ws.Range("E51").Value = "  +17.07%  "

# Reset style on cells that required a leading apostrophe to avoid numeric coercion,
# so they keep the default (unstyled) appearance matching the original workbook.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
